# Fruta / hortaliza, semanal
# Permute rows 2-20 (columns D,H,J,K,L,M,N,O,P,Q) so that each target row
# receives the values that used to belong to a different source row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (1-based worksheet row numbers)
$map = @{
    2  = 19
    3  = 14
    4  = 20
    5  = 17
    6  = 11
    7  = 13
    8  = 7
    9  = 2
    10 = 16
    11 = 18
    12 = 9
    13 = 10
    14 = 12
    15 = 4
    16 = 3
    17 = 5
    18 = 6
    19 = 15
    20 = 8
}

# Columns that change as part of the permutation (by 1-based column index)
# D=4, H=8, J=10, K=11, L=12, M=13, N=14, O=15, P=16, Q=17
$cols = @(4, 8, 10, 11, 12, 13, 14, 15, 16, 17)

# Snapshot every relevant cell's current value before writing anything,
# since sources and targets overlap.
$snapshot = @{}
foreach ($r in 2..20) {
    foreach ($c in $cols) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value()
    }
}

foreach ($target in 2..20) {
    $source = $map[$target]
    foreach ($c in $cols) {
        $ws.Cells.Item($target, $c).Value = $snapshot["$source-$c"]
    }
}
